# Populate the new column O ("vt" config for a third lidar) by mirroring
# column N's formatting, then filling in the values (most mirror column N,
# a few rows carry their own new data), and finally drop the now-unused
# column P / AA21 placeholders so the sheet's used range shrinks back down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clone N's cell formatting onto O (styles only, values come next) ---
$ws.Range("N1:N37").Copy()
$ws.Range("O1:O37").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R21").Copy()
$ws.Range("Q21").PasteSpecial(-4122)      # xlPasteFormats (Q21 loses its old "O/P-group" style)
$excel.CutCopyMode = $false

# --- 2. Values that mirror column N exactly ---
$mirrorRows = 6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,36
foreach ($r in $mirrorRows) {
    $ws.Cells.Item($r, 14).Value2 = $ws.Cells.Item($r, 14).Value2
    $ws.Cells.Item($r, 15).Value2 = $ws.Cells.Item($r, 14).Value2
}

# string-valued rows that mirror column N (regex id, name, data levels, range_name)
$mirrorTextRows = 1,4,34,35,37
foreach ($r in $mirrorTextRows) {
    $ws.Cells.Item($r, 15).Value2 = $ws.Cells.Item($r, 14).Text
}

# --- 3. Values specific to the new column (not copied from N) ---
$ws.Range("O2").Value2 = 20231023
$ws.Range("O3").Value2 = 20231110
$ws.Range("O5").Value2 = -20.5

# --- 4. Drop the now-empty column P placeholders and the stray AA21 ---
$ws.Range("P1:P37").Clear()
$ws.Range("AA21").Clear()

# --- 5. Column width: O should take on N's (wider) bestFit width, P reverts ---
$ws.Columns.Item(15).ColumnWidth = $ws.Columns.Item(14).ColumnWidth

# --- 6. View bookkeeping to match the saved selection ---
$ws.Range("M4").Select()
